$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Source:" footer label (shared string reused on B211) ---
$ws.Cells.Item(211, 2).Value = "Source: 2021-04-13"

# --- Revise existing week 11 / week 12 figures (columns D & E) ---
$ws.Cells.Item(212, 4).Value = 28748

$ws.Cells.Item(213, 4).Value = 42125
$ws.Cells.Item(213, 5).Value = 8.2

$ws.Cells.Item(214, 4).Value = 119829

$ws.Cells.Item(215, 4).Value = 172201
$ws.Cells.Item(215, 5).Value = 8.6

$ws.Cells.Item(216, 4).Value = 91935
$ws.Cells.Item(216, 5).Value = 6.8

$ws.Cells.Item(217, 4).Value = 44687
$ws.Cells.Item(217, 5).Value = 5.3

$ws.Cells.Item(219, 4).Value = 27609
$ws.Cells.Item(219, 5).Value = 6.7

$ws.Cells.Item(220, 4).Value = 47250
$ws.Cells.Item(220, 5).Value = 9.4

$ws.Cells.Item(221, 4).Value = 131500
$ws.Cells.Item(221, 5).Value = 10

$ws.Cells.Item(222, 4).Value = 187287
$ws.Cells.Item(222, 5).Value = 9.8

$ws.Cells.Item(223, 4).Value = 95351
$ws.Cells.Item(223, 5).Value = 8

$ws.Cells.Item(224, 4).Value = 46395
$ws.Cells.Item(224, 5).Value = 5.7

# --- Remove the "CAREFUL EASTER WEEK" notice and its yellow highlight ---
$ws.Cells.Item(225, 2).ClearContents()
$ws.Cells.Item(225, 2).Interior.Pattern = -4142

# --- Append the new week 13 data block ---
$newRows = @(
    @(2021, 13, "0-4",   21632, 7.1),
    @(2021, 13, "5-14",  30455, 11.1),
    @(2021, 13, "15-34", 107306, 11.6),
    @(2021, 13, "35-59", 151992, 11.5),
    @(2021, 13, "60-79", 77135, 9.7),
    @(2021, 13, ">=80",  38994, 6.1)
)

$r = 226
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# --- Move the active selection down to the new bottom of the frozen pane,
#     matching where a user would end up after appending this block ---
$ws.Range("E228").Select()
